$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (Calculus story): fill in the Won't Have / Could Have fields, and
#     rewrite the user story itself. ---
$ws.Range("H4").Value = "Fairies living within the code."
$ws.Range("G4").Value = "Error Checking on Input values to verify they are in correct format before starting processing."
$ws.Range("A4").Value = "I Sam, as a Mathematics , major have difficulties and struggle when first attempting to solve Calculus equations.This  Derivative/Integral Caluculator is meant to assist those who need to see step by step procedures in order to better understand the process which the"

# --- Row 9 (Physics story): rewrite the user story. ---
$ws.Range("A9").Value = "I Sam, as a user who must take Physics to obtain my degree have trouble because the Fundementals of Physics can be  frusterating for students who tend to be visual learners. This web application is meant to alleviate those troubles by adding graphics to the calculations to help students visualize the concepts by seeing graphics of how interactions of objects occur and the different forces/energies that they interact with as well as the objects themselves. "

# --- Row 5 gains Must Have / Should Have / Could Have entries (the old G4 text
#     moves down to F5 as the "Should Have"). ---
$ws.Range("E5").Value = "Must be integrated within a web application. "
$ws.Range("G5").Value = "Additional resource options if the student does not understand the steps."
$ws.Range("F5").Value = "Methods showing the steps that are taking place throughout the calculations."

# Row 5 grows to accommodate the wrapped text.
$ws.Rows("5:5").RowHeight = 80

# --- Final selection ---
$ws.Range("G5").Select()
